$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.766.15'
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").Value = '1.699.80'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.58%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.56'
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4026'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.001'
$ws.Range("E9").Value = '  -0.31%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.471'
$ws.Range("E10").Value = '  -0.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.73'
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.02'
$ws.Range("E13").Value = '  +4.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.520'
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.000'
$ws.Range("E15").Value = '  +0.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001345'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = '1.699.89'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.87'
$ws.Range("E18").Value = '  -2.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07183'
$ws.Range("E19").Value = '  +1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.79'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.345'
$ws.Range("E21").Value = '  +1.71%  '

$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.38'
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").Value = '24.758.82'
$ws.Range("E24").Value = '  +1.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.366'
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.920'
$ws.Range("E26").Value = '  -0.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.12'
$ws.Range("E27").Value = '  +2.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.165'
$ws.Range("E28").Value = '  +18.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.41'
$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.617'
$ws.Range("E30").Value = '  +1.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '144.16'
$ws.Range("E31").Value = '  +5.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.431'
$ws.Range("E32").Value = '  +24.63%  '

$ws.Range("D33").Value = '1.878.65'
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08639'
$ws.Range("E34").Value = '  -1.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.327'
$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03181'
$ws.Range("E36").Value = '  +10.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.040'
$ws.Range("E37").Value = '  +2.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2838'
$ws.Range("E38").Value = '  +1.45%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09435'
$ws.Range("E39").Value = '  +3.83%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.75'
$ws.Range("E40").Value = '  +0.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8279'
$ws.Range("E41").Value = '  +5.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.18'
$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("E43").Value = '  +1.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.66'
$ws.Range("E44").Value = '  +6.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.689'
$ws.Range("E45").Value = '  +4.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7419'
$ws.Range("E46").Value = '  +3.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.211'
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.387'
$ws.Range("E48").Value = '  +3.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9981'
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.08375'
$ws.Range("E50").Value = '  +4.52%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.63'
$ws.Range("E51").Value = '  +1.18%  '
